# Apply the commit changes:
# 1. Rename B1 header on "Weekly Quantity" sheet to "Weekly_PO_Qty"
# 2. Rename B1 header on "Monthly Trend" sheet to "Monthly_PO_Qty"
# 3. Add a new "PO Forecast" worksheet with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Match the bold/centered header style used on the other sheets
$wsWeekly.Range("B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45158.99999999999,204,-32.55064043525615,430.2274596387181),
    @(45165.99999999999,203,-31.57354866825045,438.6746964330746),
    @(45172.99999999999,203,-23.16733325399984,430.7941638761689),
    @(45200.99999999999,201,-46.56725261573204,427.5976145052636),
    @(45207.99999999999,200,-35.91328321814197,436.6316400844664),
    @(45221.99999999999,199,-26.73687865278275,434.4573990036204),
    @(45228.99999999999,199,-43.91352446853358,432.8260834266584),
    @(45235.99999999999,198,-46.31162974531042,411.7519603941259),
    @(45242.99999999999,197,-25.0230231751175,428.7431866337008),
    @(45249.99999999999,197,-34.40153720460361,435.2689336150852),
    @(45256.99999999999,196,-33.34865263757575,439.9420524601638),
    @(45270.99999999999,195,-44.60066707468547,434.3011085949782),
    @(45277.99999999999,195,-46.21800887216153,420.6668224382052),
    @(45298.99999999999,193,-29.84112207494589,413.1754443191113),
    @(45305.99999999999,193,-35.75156420509585,425.0086223511225),
    @(45319.99999999999,192,-27.2108542331915,418.9555306318125),
    @(45333.99999999999,191,-34.06936478958026,424.9151780099792),
    @(45340.99999999999,190,-45.9025589085552,412.1937658356766),
    @(45347.99999999999,190,-53.8234447469605,419.0005034127905),
    @(45354.99999999999,189,-48.1851564669348,420.3684944147423),
    @(45361.99999999999,189,-48.87264857988073,414.6682500799454),
    @(45375.99999999999,188,-32.73731378347163,434.117470385173),
    @(45382.99999999999,187,-40.45192464703468,420.3484965458568),
    @(45396.99999999999,186,-52.98725753684488,424.8420181510267),
    @(45410.99999999999,185,-46.70109803760162,417.6370043623736),
    @(45417.99999999999,185,-26.67641808651074,409.6060038577995),
    @(45424.99999999999,184,-54.45717238737602,410.9560910374474),
    @(45431.99999999999,184,-51.1127269201232,413.2741888904006),
    @(45445.99999999999,183,-55.45707914954777,418.6353770153341),
    @(45452.99999999999,182,-57.00614311040298,415.3883828164034),
    @(45480.99999999999,180,-48.89319715366598,413.1283559828217),
    @(45487.99999999999,180,-62.8612336779255,404.2649824417459),
    @(45494.99999999999,179,-43.71360496044021,414.0649747990662),
    @(45501.99999999999,178,-54.68764460790966,400.3957720020593),
    @(45515.99999999999,177,-40.67955229445016,411.788921877631),
    @(45543.99999999999,175,-63.18200937976218,409.6965213115236),
    @(45550.99999999999,175,-54.43267027615747,415.9905536813082),
    @(45557.99999999999,174,-58.55226845031756,411.1299913892874),
    @(45564.99999999999,174,-45.49751795997836,387.6064504158395),
    @(45578.99999999999,173,-63.34652630564956,410.0981814388686),
    @(45585.99999999999,172,-45.89084094578814,389.7733078430113),
    @(45592.99999999999,172,-64.41444638341024,391.075287280995),
    @(45599.99999999999,171,-48.94701153990947,389.6678759731209),
    @(45620.99999999999,170,-39.11343603665462,416.7372449146417),
    @(45627.99999999999,169,-61.08280828530975,411.246086906633),
    @(45634.99999999999,169,-51.46442994181664,406.0469483277885),
    @(45641.99999999999,168,-70.34364112938175,398.6860508582676),
    @(45648.99999999999,168,-72.62668243985401,388.8123680679867),
    @(45655.99999999999,167,-75.72779800613267,386.6159846176084),
    @(45662.99999999999,167,-83.52236473663471,417.002029209577),
    @(45669.99999999999,166,-75.73951742308358,400.8996939062923),
    @(45676.99999999999,166,-82.11479135944712,397.0263678307568)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$newSheet.Range("A1").Select()
